# Apply the "hacked release schedule" data edits described by the commit:
#   - On the Picarro sheet, mark BSD (col B) and TAC (col E) as "x" (no
#     data) for the ch4 row (row 10), matching the existing "x" pattern
#     already present on neighbouring rows.
#   - On the LGR sheet, mark TAC (col B) as "x" for both the n2o (row 9)
#     and co (row 10) rows.
#   - Leave the Picarro tab as the active/selected sheet (it becomes the
#     active tab after these edits), with the cursor left on D22.

$wb = $excel.ActiveWorkbook

# --- Picarro sheet -------------------------------------------------------
$picarro = $wb.Worksheets.Item("Picarro")
$picarro.Activate()

$picarro.Range("B10").Value = "x"
$picarro.Range("E10").Value = "x"

# Leave the selection where the author left it after editing.
$picarro.Range("D22").Select() | Out-Null

# --- LGR sheet -------------------------------------------------------
$lgr = $wb.Worksheets.Item("LGR")
$lgr.Range("B9").Value = "x"
$lgr.Range("B10").Value = "x"

# Re-activate Picarro so it ends up as the saved active tab.
$picarro.Activate()
